# Update "想去人数" (interest count) figures in column F, refreshed from a
# newer scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1 (tab order) = "展览"
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F2").Value = 387
$wsExhibit.Range("F3").Value = 115
$wsExhibit.Range("F5").Value = 17
$wsExhibit.Range("F6").Value = 0
$wsExhibit.Range("F8").Value = 144
$wsExhibit.Range("F10").Value = 0

# Sheet 4 (tab order) = "全部类型"
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 0
$wsAll.Range("F4").Value = 0
$wsAll.Range("F6").Value = 23
$wsAll.Range("F7").Value = 0
$wsAll.Range("F8").Value = 144
$wsAll.Range("F9").Value = 0
$wsAll.Range("F10").Value = 483
